$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3148.0396
$ws.Range("I138").Value = 1603.88
$ws.Range("J138").Value = 3904.9805
$ws.Range("K138").Value = 4811.64
$ws.Range("L138").Value = 11714.9415
$ws.Range("M138").Value = 328.3599999999997
$ws.Range("N138").Value = -21994.9415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 704.6842
$ws.Range("I2").Value = 490.75
$ws.Range("J2").Value = 1071.4286
$ws.Range("K2").Value = 490.75
$ws.Range("L2").Value = 1071.4286
$ws.Range("M2").Value = -377.75
$ws.Range("N2").Value = -1297.4286

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H45").Value = 20897454
$ws.Range("I45").Value = 25719406
$ws.Range("K45").Value = 25719406
$ws.Range("M45").Value = -25719029

$ws.Range("H62").Value = 19800
$ws.Range("J62").Value = 19800
$ws.Range("L62").Value = 19800
$ws.Range("N62").Value = -21048

$ws.Range("H63").Value = 2652
$ws.Range("I63").Value = 2652
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2652
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1966
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 18888
$ws.Range("J64").Value = 18888
$ws.Range("L64").Value = 18888
$ws.Range("N64").Value = -19384

$ws.Range("H65").Value = 19800
$ws.Range("J65").Value = 19800
$ws.Range("L65").Value = 59400
$ws.Range("N65").Value = -65640

$ws.Range("H66").Value = 2652
$ws.Range("I66").Value = 2652
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 13260
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -9828
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 18888
$ws.Range("J67").Value = 18888
$ws.Range("L67").Value = 18888
$ws.Range("N67").Value = -20604

$ws.Range("H80").Value = 35110
$ws.Range("J80").Value = 35110
$ws.Range("L80").Value = 35110
$ws.Range("N80").Value = -37106

$ws.Range("H82").Value = 47770
$ws.Range("J82").Value = 47770
$ws.Range("L82").Value = 47770
$ws.Range("N82").Value = -48492

$ws.Range("H83").Value = 35110
$ws.Range("J83").Value = 35110
$ws.Range("L83").Value = 105330
$ws.Range("N83").Value = -115314

$ws.Range("H85").Value = 47770
$ws.Range("J85").Value = 47770
$ws.Range("L85").Value = 47770
$ws.Range("N85").Value = -50266

$ws.Range("H87").Value = 44999
$ws.Range("J87").Value = 44999
$ws.Range("L87").Value = 44999
$ws.Range("N87").Value = -47495

$ws.Range("H90").Value = 44999
$ws.Range("J90").Value = 44999
$ws.Range("L90").Value = 134997
$ws.Range("N90").Value = -147477

$ws.Range("H102").Value = 1285.1
$ws.Range("I102").Value = 1285.1
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1285.1
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 336.9000000000001
$ws.Range("N102").ClearContents()

$ws.Range("H110").Value = 1337.862
$ws.Range("I110").Value = 835.56525
$ws.Range("J110").Value = 3263.3333
$ws.Range("K110").Value = 835.56525
$ws.Range("L110").Value = 3263.3333
$ws.Range("M110").Value = 1209.43475
$ws.Range("N110").Value = -7353.3333

$ws.Range("H116").Value = 704.6842
$ws.Range("I116").Value = 490.75
$ws.Range("J116").Value = 1071.4286
$ws.Range("K116").Value = 490.75
$ws.Range("L116").Value = 1071.4286
$ws.Range("M116").Value = 1803.25
$ws.Range("N116").Value = -5659.4286

$ws.Range("H122").Value = 1072.625
$ws.Range("I122").Value = 1061.5454
$ws.Range("J122").Value = 1097
$ws.Range("K122").Value = 3184.6362
$ws.Range("L122").Value = 3291
$ws.Range("M122").Value = -734.6361999999999
$ws.Range("N122").Value = -8191

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 704.6842
$ws.Range("I3").Value = 490.75
$ws.Range("J3").Value = 1071.4286
$ws.Range("K3").Value = 490.75
$ws.Range("L3").Value = 1071.4286
$ws.Range("M3").Value = -376.75
$ws.Range("N3").Value = -1299.4286

$ws.Range("H86").Value = 4655562
$ws.Range("J86").Value = 7755104
$ws.Range("L86").Value = 7755104
$ws.Range("N86").Value = -7757350

$ws.Range("H89").Value = 4655562
$ws.Range("J89").Value = 7755104
$ws.Range("L89").Value = 38775520
$ws.Range("N89").Value = -38786752

$ws.Range("H99").Value = 1862.125
$ws.Range("I99").Value = 1100
$ws.Range("J99").Value = 2175.9412
$ws.Range("K99").Value = 1100
$ws.Range("L99").Value = 2175.9412
$ws.Range("M99").Value = 398
$ws.Range("N99").Value = -5171.9412

$ws.Range("H105").Value = 40002664
$ws.Range("I105").Value = 2838
$ws.Range("K105").Value = 2838
$ws.Range("M105").Value = -1091

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8333894.5
$ws.Range("I113").Value = 6944945
$ws.Range("K113").Value = 20834835
$ws.Range("M113").Value = -20832665

$ws.Range("H129").Value = 25642198
$ws.Range("I129").Value = 430
$ws.Range("J129").Value = 27779012
$ws.Range("K129").Value = 1290
$ws.Range("L129").Value = 83337036
$ws.Range("M129").Value = 3710
$ws.Range("N129").Value = -83347036

$ws.Range("H131").Value = 707.28
$ws.Range("I131").Value = 336.66666
$ws.Range("J131").Value = 730.93616
$ws.Range("K131").Value = 1009.99998
$ws.Range("L131").Value = 2192.80848
$ws.Range("M131").Value = 4030.00002
$ws.Range("N131").Value = -12272.80848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1248.375
$ws.Range("I102").Value = 1231.6
$ws.Range("K102").Value = 1231.6
$ws.Range("M102").Value = 390.4000000000001

$ws.Range("H113").Value = 2411.125
$ws.Range("I113").Value = 2077.6667
$ws.Range("J113").Value = 2839.8572
$ws.Range("K113").Value = 2077.6667
$ws.Range("L113").Value = 2839.8572
$ws.Range("M113").Value = 92.33329999999978
$ws.Range("N113").Value = -7179.8572

$ws.Range("H126").Value = 3332
$ws.Range("I126").Value = 3693.6667
$ws.Range("J126").Value = 3151.1667
$ws.Range("K126").Value = 11081.0001
$ws.Range("L126").Value = 9453.500100000001
$ws.Range("M126").Value = -8611.000100000001
$ws.Range("N126").Value = -14393.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1249.6666
$ws.Range("I82").Value = 1104.5
$ws.Range("K82").Value = 1104.5
$ws.Range("M82").Value = -743.5

$ws.Range("H85").Value = 1249.6666
$ws.Range("I85").Value = 1104.5
$ws.Range("K85").Value = 1104.5
$ws.Range("M85").Value = 143.5

$ws.Range("H100").Value = 2673.7646
$ws.Range("I100").Value = 2109.5
$ws.Range("J100").Value = 3708.25
$ws.Range("K100").Value = 2109.5
$ws.Range("L100").Value = 3708.25
$ws.Range("M100").Value = -1568.5
$ws.Range("N100").Value = -4790.25
